{"js": "// Update the worksheet date and every multiplication problem's text to the\n// new values from the authoring diff. Each \"old\" string occurs exactly once\n// in the document, so a direct search + replace is unambiguous and safe.\nconst replacements = [\n  [\"2025-07-25 Friday\", \"2025-07-26 Saturday\"],\n  [\"92\u00d747=\", \"52\u00d793=\"],\n  [\"23\u00d765=\", \"78\u00d711=\"],\n  [\"52\u00d725=\", \"52\u00d731=\"],\n  [\"71\u00d777=\", \"42\u00d782=\"],\n  [\"67\u00d763=\", \"46\u00d749=\"],\n  [\"88\u00d718=\", \"86\u00d763=\"],\n  [\"41\u00d757=\", \"93\u00d749=\"],\n  [\"88\u00d744=\", \"39\u00d762=\"],\n  [\"33\u00d715=\", \"52\u00d772=\"],\n  [\"38\u00d794=\", \"57\u00d798=\"],\n  [\"55\u00d744=\", \"19\u00d786=\"],\n  [\"72\u00d787=\", \"55\u00d784=\"],\n  [\"67\u00d751=\", \"66\u00d738=\"],\n  [\"38\u00d777=\", \"95\u00d781=\"],\n  [\"55\u00d771=\", \"48\u00d789=\"],\n  [\"75\u00d745=\", \"28\u00d715=\"],\n  [\"44\u00d744=\", \"46\u00d726=\"],\n  [\"57\u00d714=\", \"74\u00d712=\"],\n  [\"19\u00d711=\", \"65\u00d712=\"],\n  [\"92\u00d717=\", \"31\u00d722=\"],\n  [\"64\u00d736=\", \"52\u00d718=\"],\n  [\"65\u00d746=\", \"73\u00d754=\"],\n  [\"35\u00d732=\", \"63\u00d798=\"],\n  [\"75\u00d786=\", \"87\u00d754=\"],\n  [\"28\u00d737=\", \"71\u00d779=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every multiplication problem in the table\n# with the updated values from the authoring diff. Each old string occurs\n# exactly once in the document, so Find/Replace (wdReplaceAll, scoped to a\n# single exact match) is unambiguous for every pair.\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @(\"2025-07-25 Friday\", \"2025-07-26 Saturday\"),\n    @(\"92\u00d747=\", \"52\u00d793=\"),\n    @(\"23\u00d765=\", \"78\u00d711=\"),\n    @(\"52\u00d725=\", \"52\u00d731=\"),\n    @(\"71\u00d777=\", \"42\u00d782=\"),\n    @(\"67\u00d763=\", \"46\u00d749=\"),\n    @(\"88\u00d718=\", \"86\u00d763=\"),\n    @(\"41\u00d757=\", \"93\u00d749=\"),\n    @(\"88\u00d744=\", \"39\u00d762=\"),\n    @(\"33\u00d715=\", \"52\u00d772=\"),\n    @(\"38\u00d794=\", \"57\u00d798=\"),\n    @(\"55\u00d744=\", \"19\u00d786=\"),\n    @(\"72\u00d787=\", \"55\u00d784=\"),\n    @(\"67\u00d751=\", \"66\u00d738=\"),\n    @(\"38\u00d777=\", \"95\u00d781=\"),\n    @(\"55\u00d771=\", \"48\u00d789=\"),\n    @(\"75\u00d745=\", \"28\u00d715=\"),\n    @(\"44\u00d744=\", \"46\u00d726=\"),\n    @(\"57\u00d714=\", \"74\u00d712=\"),\n    @(\"19\u00d711=\", \"65\u00d712=\"),\n    @(\"92\u00d717=\", \"31\u00d722=\"),\n    @(\"64\u00d736=\", \"52\u00d718=\"),\n    @(\"65\u00d746=\", \"73\u00d754=\"),\n    @(\"35\u00d732=\", \"63\u00d798=\"),\n    @(\"75\u00d786=\", \"87\u00d754=\"),\n    @(\"28\u00d737=\", \"71\u00d779=\")\n)\n\n$d = $word.ActiveDocument\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
